# Updated cryptos list on Thu Oct 12 06:46:22 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $b, $c, $d, $e) {
    if ($null -ne $b) { $ws.Cells.Item($row, 2).Value = $b }
    if ($null -ne $c) { $ws.Cells.Item($row, 3).Value = $c }
    if ($null -ne $d) { $ws.Cells.Item($row, 4).Value = $d }
    if ($null -ne $e) { $ws.Cells.Item($row, 5).Value = $e }
}

# Row 2 - Bitcoin
Set-Row 2 $null $null "26.915.05" "  -0.93%  "
# Row 3 - Ethereum
Set-Row 3 $null $null "1.566.10" "  +0.53%  "
# Row 4 - TetherUSD
Set-Row 4 $null $null $null "  -0.20%  "
# Row 5 - BNB
Set-Row 5 $null $null "206.16" "  -0.15%  "
# Row 6 - XRP
Set-Row 6 $null $null $null "  -0.44%  "
# Row 7 - USDC
Set-Row 7 $null $null $null "  -0.19%  "
# Row 8 - Solana
Set-Row 8 $null $null $null "  -1.47%  "
# Row 9 - Cardano
Set-Row 9 $null $null $null "  -0.25%  "
# Row 10 - Dogecoin
Set-Row 10 $null $null $null "  -1.02%  "
# Row 11 - TRON
Set-Row 11 $null $null $null "  +0.13%  "
# Row 12 - WrappedliquidstakedEther2.0
Set-Row 12 $null $null "1.788.99" "  +0.55%  "
# Row 13 - WrappedEther
Set-Row 13 $null $null "1.564.85" "  +0.50%  "
# Row 14 - Polkadot
Set-Row 14 $null $null $null "  -0.84%  "
# Row 15 - Polygon
Set-Row 15 $null $null $null "  +0.08%  "
# Row 16 - WrappedBTC
Set-Row 16 $null $null "26.903.93" "  -0.90%  "
# Row 17 - Litecoin
Set-Row 17 $null $null "61.37" "  -2.36%  "
# Row 18 - BitcoinCash
Set-Row 18 $null $null "215.85" "  -0.22%  "
# Row 19 - Chainlink
Set-Row 19 $null $null "7.37" "  +1.91%  "
# Row 20 - ShibaInu
Set-Row 20 $null $null $null "  -0.45%  "
# Row 21 - Dai
Set-Row 21 $null $null $null "  -0.19%  "
# Row 22 - Uniswap
Set-Row 22 $null $null $null "  +0.90%  "
# Row 23 - Avalanche
Set-Row 23 $null $null $null "  -1.38%  "
# Row 24 - Toncoin
Set-Row 24 $null $null $null "  +0.84%  "
# Row 25 - Monero
Set-Row 25 $null $null "153.98" "  +1.35%  "
# Row 26 - Cosmos
Set-Row 26 $null $null $null "  +1.96%  "
# Row 27 - EthereumClassic
Set-Row 27 $null $null $null "  +0.41%  "
# Row 28 - BinanceUSD
Set-Row 28 $null $null $null "  -0.21%  "
# Row 29 - Stellar
Set-Row 29 $null $null $null "  -0.66%  "
# Row 30 - Hedera
Set-Row 30 $null $null "0.0466" "  +0.94%  "
# Row 31 - PancakeSwap
Set-Row 31 $null $null "1.11" "  -3.09%  "
# Row 32 - Filecoin
Set-Row 32 $null $null $null "  -0.24%  "
# Row 33 - Maker
Set-Row 33 $null $null "1.400.28" "  +1.64%  "
# Row 34 - InternetComputer(DFINITY)
Set-Row 34 $null $null $null "  +0.09%  "
# Row 35 - LidoDAOToken
Set-Row 35 $null $null $null "  -0.62%  "
# Row 36 - HuobiToken
Set-Row 36 $null $null $null "  -0.32%  "
# Row 37 - TrustWalletToken
Set-Row 37 $null $null "0.921" "  -2.79%  "
# Row 38 - VeChain
Set-Row 38 $null $null $null "  -0.53%  "
# Row 39 - ImmutableX
Set-Row 39 $null $null "0.531" "  +3.12%  "
# Row 40 - ARBITRUM
Set-Row 40 $null $null $null "  +0.33%  "
# Row 41 - PaxDollar
Set-Row 41 $null $null $null "  -0.18%  "
# Row 42 - WEMIXToken
Set-Row 42 $null $null $null "  +0.82%  "
# Row 43 - FraxShare
Set-Row 43 $null $null $null "  +6.49%  "
# Row 44 - RenderToken
Set-Row 44 $null $null $null "  -0.01%  "

# Row 45 - Aave -> MXToken (rows 45/46 swap ranking order)
Set-Row 45 "MXToken" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" "2.18" "  +1.09%  "
# Row 46 - MXToken -> Aave
Set-Row 46 "Aave" "https://coinranking.com/coin/ixgUfzmLR+aave-aave" "63.91" "  +1.07%  "

# Row 47 - RocketPoolETH
Set-Row 47 $null $null "1.702.48" "  +0.70%  "
# Row 48 - Quant
Set-Row 48 $null $null "86.87" "  +1.72%  "

# Row 49 - Cronos -> BabyDogeCoin (rows 49/50 swap ranking order)
Set-Row 49 "BabyDogeCoin" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge" "0.0₇0982" "  +0.10%  "
# Row 50 - BabyDogeCoin -> Cronos
Set-Row 50 "Cronos" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro" "0.0502" "  +2.10%  "

# Row 51 - Algorand
Set-Row 51 $null $null "0.0955" "  +1.43%  "
